$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the newly-added evaluation numbers for rows 25-28 (columns B:G)
$ws.Range("B25").Value = 0.69079999999999997
$ws.Range("C25").Value = 0.52659999999999996
$ws.Range("D25").Value = 0.69120000000000004
$ws.Range("E25").Value = 0.52510000000000001
$ws.Range("F25").Value = 0.69089999999999996
$ws.Range("G25").Value = 0.5262

$ws.Range("B26").Value = 0.69320000000000004
$ws.Range("C26").Value = 0.50160000000000005
$ws.Range("D26").Value = 0.69289999999999996
$ws.Range("E26").Value = 0.50700000000000001
$ws.Range("F26").Value = 0.69330000000000003
$ws.Range("G26").Value = 0.50009999999999999

$ws.Range("B27").Value = 0.69920000000000004
$ws.Range("C27").Value = 0.53359999999999996
$ws.Range("D27").Value = 0.6905
$ws.Range("E27").Value = 0.5333
$ws.Range("F27").Value = 0.69059999999999999
$ws.Range("G27").Value = 0.52859999999999996

$ws.Range("B28").Value = 0.69289999999999996
$ws.Range("C28").Value = 0.51039999999999996
$ws.Range("D28").Value = 0.69240000000000002
$ws.Range("E28").Value = 0.52029999999999998
$ws.Range("F28").Value = 0.69299999999999995
$ws.Range("G28").Value = 0.50739999999999996

# Update the active selection on the sheet to match the saved view state
$ws.Range("C29").Select()
